# Update countries & provincias Spain
# The underlying COVID dataset was refreshed. This produced:
#  - updated case numbers for several countries (Espana, Libano, Albania, Brunei, Georgia)
#  - a newly-tracked country (Eslovenia) that is re-sorted into the
#    descending "Casos totales" ranking, which cascades row shifts for the
#    rows that used to sit below it (Estonia, Argelia)
#  - the same cascading shift for Kuwait (now ranked above Taiwan /
#    Republica de Chipre) and for Mayotte (now ranked above Camboya)
#  - the "last updated" timestamp footer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($rowIndex, $country, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($rowIndex, 1).Value = $country
    $ws.Cells.Item($rowIndex, 2).Value = $casosTotales
    $ws.Cells.Item($rowIndex, 3).Value = $nuevosCasos
    $ws.Cells.Item($rowIndex, 4).Value = $casosActivos
    $ws.Cells.Item($rowIndex, 5).Value = $recuperados
    $ws.Cells.Item($rowIndex, 6).Value = $casosCriticos
    $ws.Cells.Item($rowIndex, 7).Value = $muertesHoy
    $ws.Cells.Item($rowIndex, 8).Value = $muertes
}

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 11:50"

# Spain refreshed numbers
Set-Row 6 "España" 110238 6120 26743 73492 5872 616 10003

# Eslovenia newly ranked ahead of Estonia / Argelia (cascading shift)
Set-Row 55 "Eslovenia" 897 56 10 871 29 1 16
Set-Row 56 "Estonia" 858 79 45 802 16 6 11
Set-Row 57 "Argelia" 847 0 61 728 0 0 58

# Libano refreshed numbers
Set-Row 72 "Libano" 494 15 43 435 3 2 16

# Kuwait newly ranked ahead of Taiwan / Republica de Chipre (cascading shift)
Set-Row 84 "Kuwait" 342 25 81 261 15 0 0
Set-Row 85 "Taiwan" 339 10 50 284 0 0 5
Set-Row 86 "Republica de Chipre" 320 0 28 283 11 0 9

# Albania refreshed numbers
Set-Row 90 "Albania" 277 18 67 194 7 1 16

# Brunei / Georgia refreshed numbers
Set-Row 112 "Brunei" 133 2 56 76 3 0 1
Set-Row 113 "Georgia" 130 13 26 104 6 0 0

# Mayotte newly ranked ahead of Camboya (cascading shift)
Set-Row 118 "Mayotte" 116 15 10 105 3 0 1
Set-Row 119 "Camboya" 110 1 34 76 1 0 0
